$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the wrapped Times-New-Roman format already used in column D (rows 5-8)
# onto D9 before writing its value, so it reuses the existing style instead of
# minting a new cellXf.
$ws.Range("D8").Copy()
$ws.Range("D9").PasteSpecial(-4122)

$mu = [char]0x03BC

$ws.Range("D9").Value = "8.0725513 pixels per " + $mu + "m"
$ws.Range("E9").Value = "-p [1500,700,3,250]"
$ws.Range("F9").Value = 0.55
$ws.Range("G9").Value = 0.1
$ws.Range("I9").Value = "y"
$ws.Range("J9").Value = "Almost all bad due to debris and angle, two might be okay, but area seems overestimated on all objects, and neither avicularia on autozooid 1 are detected, even with confidence 0, strictness 1.00 or autofilter turned off."
$ws.Range("K9").Value = 4
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 1
$ws.Range("N9").Value = 1

$ws.Rows.Item(9).RowHeight = 23.85
